$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.430.26'
$ws.Range('E2').Value = '  -0.82%  '
$ws.Range('D3').Value = '2.331.27'
$ws.Range('E3').Value = '  -0.55%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'512.37"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.28%  '
$ws.Range('D6').Value = "'132.22"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.47%  '
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('E8').Value = '  -0.76%  '
$ws.Range('E9').Value = '  -3.14%  '
$ws.Range('E10').Value = '  -0.48%  '
$ws.Range('E11').Value = '  +1.94%  '
$ws.Range('E12').Value = '  -0.53%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '2.746.11'
$ws.Range('E13').Value = '  -0.36%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').Value = "'23.54"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.46%  '
$ws.Range('D15').Value = '56.424.85'
$ws.Range('E15').Value = '  -0.60%  '
$ws.Range('D16').Value = "'0.0000133"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.14%  '
$ws.Range('D17').Value = '2.331.26'
$ws.Range('E17').Value = '  +0.74%  '
$ws.Range('D18').Value = "'10.46"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.53%  '
$ws.Range('D19').Value = "'324.23"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.83%  '
$ws.Range('D21').Value = "'6.65"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.92%  '
$ws.Range('D22').Value = "'0.996"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.38%  '
$ws.Range('D23').Value = "'61.69"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.63%  '
$ws.Range('D24').Value = "'8.67"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +11.66%  '
$ws.Range('E25').Value = '  +0.93%  '
$ws.Range('D26').Value = "'0.999"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.21%  '
$ws.Range('E27').Value = '  +6.52%  '
$ws.Range('D28').Value = "'167.55"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.48%  '
$ws.Range('B29').Value = 'PancakeSwap'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D29').Value = "'1.67"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.55%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').Value = '0.0₃0720'
$ws.Range('E30').Value = '  -2.94%  '
$ws.Range('D31').Value = "'6.10"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.40%  '
$ws.Range('D32').Value = "'18.29"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.38%  '
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('D34').Value = "'0.998"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.19%  '
$ws.Range('D36').Value = "'3.94"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.83%  '
$ws.Range('D37').Value = "'0.887"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.33%  '
$ws.Range('D38').Value = "'38.45"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.89%  '
$ws.Range('E39').Value = '  +1.01%  '
$ws.Range('D40').Value = "'150.85"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +10.28%  '
$ws.Range('D41').Value = "'0.375"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.75%  '
$ws.Range('E42').Value = '  -0.17%  '
$ws.Range('D43').Value = "'279.44"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.85%  '
$ws.Range('D44').Value = "'5.05"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.65%  '
$ws.Range('E45').Value = '  -0.87%  '
$ws.Range('E46').Value = '  -1.46%  '
$ws.Range('D47').Value = "'0.558"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.50%  '
$ws.Range('D48').Value = "'18.04"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +5.30%  '
$ws.Range('D49').Value = "'0.380"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.43%  '
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('D51').Value = "'17.19"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.36%  '
